$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "'257.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'0.05%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'26.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.52%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'4.567"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-12.66%"
$ws.Range("E4").Style = "Normal"
$ws.Range("E5").Value = "'-0.43%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.634"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-0.87%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.8583"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-1.32%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.9268"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-12.28%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.1408"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.68%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.03807"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'4.85%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07084"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-1.58%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.03172"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.70%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09144"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.84%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.001536"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.46%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'TigerCash"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.005988"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.72%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'LEO"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'3.517"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.02%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'GateToken"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'3.198"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.03%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'BTSEToken"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'2.204"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.23%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'One"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'0.01045"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'1,625.44%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-1.40%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.1293"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'0.07%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.895"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'9.70%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04213"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.65%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.001218"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-0.17%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-5.13%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-0.22%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-22.24%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D40").Value = "'0.03844"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.59%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006264"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'15.00%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1103"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.22%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002198"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-4.55%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01147"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'15.68%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005458"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.28%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.22%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.05995"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-45.13%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.1396"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'6,422.05%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.22%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.22%"
$ws.Range("E50").Style = "Normal"
